$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "vat" column header
$ws.Range("L1").Value = "vat"

# Add vat values for each product row
$ws.Range("L2").Value = 5000
$ws.Range("L3").Value = 2000
$ws.Range("L4").Value = 1000
$ws.Range("L5").Value = 3000

# Update selection to match final state (L3 selected)
$ws.Range("L3").Select()
